$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from existing header cell (H1) to new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Add new column headers
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill data rows 2-24: I = 1 (constant), J = same value as column H
for ($r = 2; $r -le 24; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
